$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_1", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_2", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_3", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_1", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_2", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_3", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_1", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_2", "Dummy"),
    @("github-gk-aks/secondgithubrepo", "prod", "PROD_ENV_3", "Dummy"),
    @("github-gk-aks/firstgithubrepo", "dev", "DEV_ENV_1", "Dummy"),
    @("github-gk-aks/firstgithubrepo", "dev", "DEV_ENV_2", "Dummy"),
    @("github-gk-aks/firstgithubrepo", "dev", "DEV_ENV_3", "Dummy"),
    @("github-gk-aks/firstgithubrepo", "staging", "STAGE_ENV_1", "Dummy"),
    @("github-gk-aks/fourthgithubrepo", "uat", "UAT_ENV_1", "Dummy"),
    @("github-gk-aks/fourthgithubrepo", "uat", "UAT_ENV_2", "Dummy"),
    @("github-gk-aks/fourthgithubrepo", "uat", "UAT_ENV_3", "Dummy"),
    @("github-gk-aks/fourthgithubrepo", "uat", "UAT_ENV_4", "Dummy")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
